# Update the "想去人数" (F column) figures on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets to the newly scraped counts.
#
# Each sheet keeps its own independent row->value mapping because the
# two sheets were scraped at (very slightly) different moments, so the
# new F45 values differ by one (1324 vs 1325) between the two sheets,
# and the "全部类型" sheet's F18 had already been refreshed earlier and
# therefore does not need to change here.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)   # 展览
$sheet4 = $wb.Worksheets.Item(4)   # 全部类型

$sheet1Updates = @{
    2  = 286
    3  = 1456
    10 = 147
    12 = 4829
    14 = 7113
    16 = 62
    18 = 586
    19 = 60
    21 = 4206
    22 = 1425
    23 = 88
    25 = 2779
    28 = 182
    29 = 405
    30 = 399
    32 = 255
    33 = 58
    35 = 1089
    37 = 1022
    38 = 93
    42 = 13
    43 = 31
    45 = 1324
    46 = 663
    47 = 29
}

$sheet4Updates = @{
    2  = 286
    3  = 1456
    10 = 147
    12 = 4829
    14 = 7113
    16 = 62
    19 = 60
    21 = 4206
    22 = 1425
    23 = 88
    25 = 2779
    28 = 182
    29 = 405
    30 = 399
    32 = 255
    33 = 58
    35 = 1089
    37 = 1022
    38 = 93
    42 = 13
    43 = 31
    45 = 1325
    46 = 663
    47 = 29
}

foreach ($row in $sheet1Updates.Keys) {
    $sheet1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

foreach ($row in $sheet4Updates.Keys) {
    $sheet4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
